$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.689.61"
$ws.Range("E2").Value = "  +0.82%  "

# Row 3
$ws.Range("D3").Value = "1.602.60"
$ws.Range("E3").Value = "  -0.08%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").Value = "'212.27"
$ws.Range("E5").Value = "  -0.35%  "

# Row 6
$ws.Range("E6").Value = "  -0.11%  "

# Row 7
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").Value = "'28.87"
$ws.Range("E8").Value = "  +8.22%  "

# Row 9
$ws.Range("E9").Value = "  +2.69%  "

# Row 10
$ws.Range("E10").Value = "  +1.45%  "

# Row 11
$ws.Range("D11").Value = "'0.0907"
$ws.Range("E11").Value = "  -0.49%  "

# Row 12
$ws.Range("D12").Value = "1.831.86"
$ws.Range("E12").Value = "  -0.13%  "

# Row 13
$ws.Range("D13").Value = "1.611.85"
$ws.Range("E13").Value = "  +0.66%  "

# Row 14
$ws.Range("D14").Value = "'0.555"
$ws.Range("E14").Value = "  +3.38%  "

# Row 15
$ws.Range("D15").Value = "29.676.65"
$ws.Range("E15").Value = "  +0.54%  "

# Row 16
$ws.Range("E16").Value = "  +1.02%  "

# Row 17
$ws.Range("D17").Value = "'64.14"
$ws.Range("E17").Value = "  +1.15%  "

# Row 18
$ws.Range("D18").Value = "'240.82"
$ws.Range("E18").Value = "  +0.58%  "

# Row 19
$ws.Range("D19").Value = "'8.08"
$ws.Range("E19").Value = "  +6.68%  "

# Row 20
$ws.Range("E20").Value = "  +1.22%  "

# Row 21
$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  +0.08%  "

# Row 22
$ws.Range("D22").Value = "'4.04"
$ws.Range("E22").Value = "  +0.27%  "

# Row 23
$ws.Range("D23").Value = "'9.50"
$ws.Range("E23").Value = "  +3.36%  "

# Row 24
$ws.Range("E24").Value = "  +2.33%  "

# Row 25
$ws.Range("D25").Value = "'155.98"
$ws.Range("E25").Value = "  +0.93%  "

# Row 26
$ws.Range("D26").Value = "'15.46"
$ws.Range("E26").Value = "  +1.37%  "

# Row 27
$ws.Range("E27").Value = "  +1.15%  "

# Row 28
$ws.Range("E28").Value = "  +1.87%  "

# Row 29
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.03%  "

# Row 30
$ws.Range("E30").Value = "  +2.08%  "

# Row 31
$ws.Range("E31").Value = "  -0.26%  "

# Row 32
$ws.Range("E32").Value = "  +0.03%  "

# Row 33
$ws.Range("D33").Value = "'3.16"
$ws.Range("E33").Value = "  +2.27%  "

# Row 34
$ws.Range("D34").Value = "1.424.67"
$ws.Range("E34").Value = "  -0.40%  "

# Row 35
$ws.Range("E35").Value = "  +3.96%  "

# Row 36
$ws.Range("E36").Value = "  -0.24%  "

# Row 37
$ws.Range("E37").Value = "  +1.89%  "

# Row 38
$ws.Range("E38").Value = "  +0.61%  "

# Row 39
$ws.Range("E39").Value = "  +2.45%  "

# Row 40
$ws.Range("D40").Value = "'0.549"
$ws.Range("E40").Value = "  +3.40%  "

# Row 41
$ws.Range("D41").Value = "'55.32"
$ws.Range("E41").Value = "  +3.38%  "

# Row 42
$ws.Range("E42").Value = "  +4.74%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'1.96"
$ws.Range("E43").Value = "  +1.35%  "

# Row 44
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "'0.820"
$ws.Range("E44").Value = "  +3.35%  "

# Row 45
$ws.Range("D45").Value = "'0.998"
$ws.Range("E45").Value = "  +0.04%  "

# Row 46
$ws.Range("D46").Value = "'67.45"
$ws.Range("E46").Value = "  +2.66%  "

# Row 47
$ws.Range("D47").Value = "'0.995"
$ws.Range("E47").Value = "  +18.77%  "

# Row 48
$ws.Range("D48").Value = "'5.41"
$ws.Range("E48").Value = "  +2.89%  "

# Row 49
$ws.Range("D49").Value = "1.740.98"
$ws.Range("E49").Value = "  -0.22%  "

# Row 50
$ws.Range("B50").Value = "mCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D50").Value = "'2.12"
$ws.Range("E50").Value = "  -1.00%  "

# Row 51
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "'86.71"
$ws.Range("E51").Value = "  -0.01%  "
